$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'30.547.38"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.14%  "
$ws.Range("D3").Value = "'1.855.11"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.39%  "
$ws.Range("D4").Value = "'1.000"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").Value = "'233.62"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.34%  "
$ws.Range("E6").Value = "  -0.05%  "
$ws.Range("D7").Value = "'0.4741"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.72%  "
$ws.Range("D8").Value = "'0.2752"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.43%  "
$ws.Range("D9").Value = "'0.06318"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.82%  "
$ws.Range("E10").Value = "  +9.67%  "
$ws.Range("D11").Value = "'1.908.61"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +3.26%  "
$ws.Range("D12").Value = "'0.07450"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.37%  "
$ws.Range("D13").Value = "'4.996"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.41%  "
$ws.Range("D14").Value = "'84.63"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.36%  "
$ws.Range("D15").Value = "'0.6266"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.13%  "
$ws.Range("D16").Value = "'30.510.63"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.18%  "
$ws.Range("D17").Value = "'245.11"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +7.58%  "
$ws.Range("D18").Value = "'0.9999"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.04%  "
$ws.Range("D19").Value = "'12.72"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.15%  "
$ws.Range("D20").Value = "'0.000007339"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.27%  "
$ws.Range("D21").Value = "'1.002"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.21%  "
$ws.Range("D22").Value = "'4.954"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.73%  "
$ws.Range("E23").Value = "  +0.63%  "
$ws.Range("D24").Value = "'9.141"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.63%  "
$ws.Range("D25").Value = "'162.96"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.33%  "
$ws.Range("D26").Value = "'18.00"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.84%  "
$ws.Range("D27").Value = "'1.879"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.69%  "
$ws.Range("D28").Value = "'0.1019"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.45%  "
$ws.Range("D29").Value = "'1.359"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.55%  "
$ws.Range("D30").Value = "'4.011"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.34%  "
$ws.Range("D31").Value = "'3.839"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.86%  "
$ws.Range("E32").Value = "  -0.82%  "
$ws.Range("D33").Value = "'1.138"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.97%  "
$ws.Range("D34").Value = "'0.7037"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.37%  "
$ws.Range("E35").Value = "  -0.27%  "
$ws.Range("D36").Value = "'0.01900"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.87%  "
$ws.Range("D37").Value = "'2.687"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.16%  "
$ws.Range("D38").Value = "'2.010"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +3.35%  "
$ws.Range("D39").Value = "'0.8769"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -3.03%  "
$ws.Range("D40").Value = "'106.86"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.03%  "
$ws.Range("D41").Value = "'0.9999"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.19%  "
$ws.Range("D42").Value = "'5.550"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.25%  "
$ws.Range("E43").Value = "  -0.29%  "
$ws.Range("D44").Value = "'7.201"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.34%  "
$ws.Range("D45").Value = "'62.88"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +4.74%  "
$ws.Range("E46").Value = "  +1.84%  "
$ws.Range("D47").Value = "'33.67"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.77%  "
$ws.Range("D48").Value = "'8.570"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.66%  "
$ws.Range("D49").Value = "'0.05536"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.43%  "
$ws.Range("D50").Value = "'1.352"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.24%  "
$ws.Range("D51").Value = "'0.3692"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.74%  "
